$d = $word.ActiveDocument

# 1. Merge the title runs "E" + "xpectation of " + "O" + "riginality" into one run.
$d.Content.Find.Execute("Expectation of Originality", $true, $false, $false, $false, $false, $true, 1, $false, "Expectation of Originality", 2) | Out-Null

# 2. Remove the existing _GoBack bookmark around "COMP-5541-2184-DD".
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 3. Add team member IDs into the empty "ID:" cells of the roster table.
$t = $d.Tables.Item(1)

# Row 3 -> Danny Shash
$t.Cell(3, 1).Range.Text = "29548912"

# Row 4 -> Xinjie Zeng: re-create the _GoBack bookmark on the (still empty)
# paragraph first, then type the ID text in front of it via the same live
# Range object so the bookmark ends up right after the new run, matching
# the target layout:
#   <w:r><w:t>27238223</w:t></w:r><w:bookmarkStart .../><w:bookmarkEnd .../>
$idRange = $t.Cell(4, 1).Range
$idRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $idRange)
$idRange.InsertBefore("27238223")

# Row 5 -> Siming Huang
$t.Cell(5, 1).Range.Text = "40081588"

# Row 6 -> Tony Lac
$t.Cell(6, 1).Range.Text = "40049123"
